$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 766.6667
$ws.Range("I29").Value = 300
$ws.Range("K29").Value = 900
$ws.Range("M29").Value = -619

$ws.Range("H33").Value = 479.0909
$ws.Range("I33").Value = 358.75
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 358.75
$ws.Range("L33").Value = 800
$ws.Range("M33").Value = -129.75
$ws.Range("N33").Value = -1258

$ws.Range("H58").Value = 97.181816
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H86").Value = 2534528.5
$ws.Range("I86").Value = 10881.1
$ws.Range("J86").Value = 15152765
$ws.Range("K86").Value = 10881.1
$ws.Range("L86").Value = 15152765
$ws.Range("M86").Value = -9758.1
$ws.Range("N86").Value = -15155011

$ws.Range("H89").Value = 2534528.5
$ws.Range("I89").Value = 10881.1
$ws.Range("J89").Value = 15152765
$ws.Range("K89").Value = 54405.5
$ws.Range("L89").Value = 75763825
$ws.Range("M89").Value = -48789.5
$ws.Range("N89").Value = -75775057

$ws.Range("H112").Value = 18594.277
$ws.Range("I112").Value = 480
$ws.Range("J112").Value = 20750.738
$ws.Range("K112").Value = 1440
$ws.Range("L112").Value = 62252.21400000001
$ws.Range("M112").Value = -332
$ws.Range("N112").Value = -64468.21400000001

$ws.Range("H137").Value = 1471.9474
$ws.Range("I137").Value = 1154.6207
$ws.Range("J137").Value = 2494.4443
$ws.Range("K137").Value = 3463.8621
$ws.Range("L137").Value = 7483.3329
$ws.Range("M137").Value = -913.8620999999998
$ws.Range("N137").Value = -12583.3329

$ws.Range("H138").Value = 3317.4
$ws.Range("I138").Value = 1805.8889
$ws.Range("J138").Value = 3876.4521
$ws.Range("K138").Value = 5417.6667
$ws.Range("L138").Value = 11629.3563
$ws.Range("M138").Value = -277.6666999999998
$ws.Range("N138").Value = -21909.3563

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13004677
$ws.Range("I32").Value = 16141252
$ws.Range("K32").Value = 16141252
$ws.Range("M32").Value = -16140965

$ws.Range("H74").Value = 1496.0322
$ws.Range("I74").Value = 1569.7778
$ws.Range("J74").Value = 1300.8235
$ws.Range("K74").Value = 1569.7778
$ws.Range("L74").Value = 1300.8235
$ws.Range("M74").Value = -695.7778000000001
$ws.Range("N74").Value = -3048.8235

$ws.Range("H77").Value = 1496.0322
$ws.Range("I77").Value = 1569.7778
$ws.Range("J77").Value = 1300.8235
$ws.Range("K77").Value = 7848.889
$ws.Range("L77").Value = 6504.1175
$ws.Range("M77").Value = -3480.889
$ws.Range("N77").Value = -15240.1175

$ws.Range("H97").Value = 1031.4546
$ws.Range("J97").Value = 1287.5
$ws.Range("L97").Value = 1287.5
$ws.Range("N97").Value = -2279.5

$ws.Range("H132").Value = 15875470
$ws.Range("I132").Value = 19231418
$ws.Range("J132").Value = 10987.182
$ws.Range("K132").Value = 57694254
$ws.Range("L132").Value = 32961.546
$ws.Range("M132").Value = -57691724
$ws.Range("N132").Value = -38021.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6074.12
$ws.Range("I80").Value = 2650.4
$ws.Range("K80").Value = 2650.4
$ws.Range("M80").Value = -1652.4

$ws.Range("H83").Value = 6074.12
$ws.Range("I83").Value = 2650.4
$ws.Range("K83").Value = 13252
$ws.Range("M83").Value = -8260

$ws.Range("H134").Value = 1485812
$ws.Range("I134").Value = 3200.389
$ws.Range("J134").Value = 5298242
$ws.Range("K134").Value = 9601.167000000001
$ws.Range("L134").Value = 15894726
$ws.Range("M134").Value = -7066.167000000001
$ws.Range("N134").Value = -15899796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H7").Value = 19457.77
$ws.Range("I7").Value = 85.85714
$ws.Range("K7").Value = 85.85714
$ws.Range("M7").Value = 27.14286

$ws.Range("H132").Value = 17550920
$ws.Range("I132").Value = 1145.1666
$ws.Range("J132").Value = 47636250
$ws.Range("K132").Value = 3435.4998
$ws.Range("L132").Value = 142908750
$ws.Range("M132").Value = -905.4998000000001
$ws.Range("N132").Value = -142913810

$ws.Range("H134").Value = 1620.1818
$ws.Range("I134").Value = 1532.1351
$ws.Range("J134").Value = 2085.5715
$ws.Range("K134").Value = 4596.4053
$ws.Range("L134").Value = 6256.7145
$ws.Range("M134").Value = -2061.4053
$ws.Range("N134").Value = -11326.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15557455
$ws.Range("I5").Value = 4167019.5
$ws.Range("J5").Value = 28575096
$ws.Range("K5").Value = 12501058.5
$ws.Range("L5").Value = 85725288
$ws.Range("M5").Value = -12500946.5
$ws.Range("N5").Value = -85725512

$ws.Range("H135").Value = 15557455
$ws.Range("I135").Value = 4167019.5
$ws.Range("J135").Value = 28575096
$ws.Range("K135").Value = 37503175.5
$ws.Range("L135").Value = 257175864
$ws.Range("M135").Value = -37500640.5
$ws.Range("N135").Value = -257180934

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7146557.5
$ws.Range("I80").Value = 4108.5835
$ws.Range("J80").Value = 50001250
$ws.Range("K80").Value = 4108.5835
$ws.Range("L80").Value = 50001250
$ws.Range("M80").Value = -3110.5835
$ws.Range("N80").Value = -50003246

$ws.Range("H83").Value = 7146557.5
$ws.Range("I83").Value = 4108.5835
$ws.Range("J83").Value = 50001250
$ws.Range("K83").Value = 20542.9175
$ws.Range("L83").Value = 250006250
$ws.Range("M83").Value = -15550.9175
$ws.Range("N83").Value = -250016234

$ws.Range("H97").Value = 523.9231
$ws.Range("I97").Value = 504.88
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 504.88
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -8.879999999999995
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5466.609
$ws.Range("I136").Value = 3986.35
$ws.Range("J136").Value = 15335
$ws.Range("K136").Value = 11959.05
$ws.Range("L136").Value = 46005
$ws.Range("M136").Value = -9409.049999999999
$ws.Range("N136").Value = -51105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11642225
$ws.Range("I132").Value = 18538466
$ws.Range("J132").Value = 4818.375
$ws.Range("K132").Value = 55615398
$ws.Range("L132").Value = 14455.125
$ws.Range("M132").Value = -55612868
$ws.Range("N132").Value = -19515.125

$ws.Range("H135").Value = 41500
$ws.Range("J135").Value = 41500
$ws.Range("L135").Value = 41500
$ws.Range("N135").Value = -51640

$ws.Range("H136").Value = 5019.62
$ws.Range("I136").Value = 10181.1
$ws.Range("J136").Value = 1578.6333
$ws.Range("K136").Value = 30543.3
$ws.Range("L136").Value = 4735.8999
$ws.Range("M136").Value = -27993.3
$ws.Range("N136").Value = -9835.8999
